# edit.ps1 - apply the scenario.docx revision described by the commit:
#   "some changes. Now one of the secondary missions is to find lost gold."
#
# Three textual edits plus one relocated "_GoBack" bookmark (the position
# Word remembers as the last edit point / cursor position before save).

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2 (used throughout below)

# --- Part III (the soldiers burn the village) -----------------------------
# "...poprzez spalenie wsi i jej mieszkańców. Oddział złożony z 50 żołnierzy
#  uzbrojonych w łuki z płonącymi strzałami stanęło przed wsią..."
# becomes
# "...poprzez spalenie wraz z mieszkańcami. Oddział złożony z 50 żołnierzy
#  uzbrojonych w łuki stanął przed wsią..."
$d.Content.Find.Execute("wsi i jej mieszkańców", $true, $false, $false, $false, $false, `
    $true, 1, $false, "wraz z mieszkańcami", 2) | Out-Null

$d.Content.Find.Execute("z płonącymi strzałami stanęło", $true, $false, $false, $false, $false, `
    $true, 1, $false, "stanął", 2) | Out-Null

# --- Part II (Hermera's side quest is changed from chickens to gold) -----
# "...za znalezienie kur, które mu uciekły z kurnika. Ktoś znowu mu..."
# becomes
# "...za znalezienie skradzionego złota. Ktoś znowu mu..."
$d.Content.Find.Execute("kur, które mu uciekły z kurnika", $true, $false, $false, $false, $false, `
    $true, 1, $false, "skradzionego złota", 2) | Out-Null

# --- Relocate the "_GoBack" bookmark (Word's "last edit location" marker) -
# It used to sit between "zostać " and "wyzwoleni." near the end of the
# document; the new save point is inside the "(akcja dzieje się ...)"
# sentence, right before the closing "i)".
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$marker = "(akcja dzieje się gdy gracz wykona wszystkie misje poboczne lub gdy wszystkie odrzuci)"
$found = $d.Content.Find.Execute($marker)
if ($d.Content.Find.Found) {
    # Content.Find leaves the located text selected in $d.Content; collapse
    # to just before the final "i)" of "odrzuci)".
    $endRange = $d.Content
    $pos = $endRange.End - 2
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}
